$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'29.497.25"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +0.04%  "

$ws.Range("D3").Value = "'1.909.96"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  -0.17%  "

$ws.Range("D4").Value = "'1.006"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  +0.58%  "

$ws.Range("D5").Value = "'326.07"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.45%  "

$ws.Range("E6").Value = "  +0.41%  "

$ws.Range("D7").Value = "'0.4841"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +1.73%  "

$ws.Range("D8").Value = "'0.4071"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -0.49%  "

$ws.Range("D9").Value = "'0.08144"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +1.46%  "

$ws.Range("E10").Value = "  +0.08%  "

$ws.Range("D11").Value = "'23.44"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +4.63%  "

$ws.Range("D12").Value = "'1.913.86"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -0.52%  "

$ws.Range("D13").Value = "'6.026"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +1.20%  "

$ws.Range("D14").Value = "'7.107"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -0.85%  "

$ws.Range("D15").Value = "'90.41"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +1.04%  "

$ws.Range("D16").Value = "'0.06796"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +2.87%  "

$ws.Range("D17").Value = "'1.006"
$ws.Range("D17").ClearFormats()

$ws.Range("E18").Value = "  +0.76%  "

$ws.Range("D19").Value = "'17.71"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -0.37%  "

$ws.Range("D20").Value = "'1.005"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +0.45%  "

$ws.Range("D21").Value = "'29.510.87"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +0.02%  "

$ws.Range("D22").Value = "'5.608"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +1.08%  "

$ws.Range("D23").Value = "'11.79"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +2.28%  "

$ws.Range("D24").Value = "'2.166"
$ws.Range("D24").ClearFormats()

$ws.Range("D25").Value = "'2.150.39"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +0.42%  "

$ws.Range("D26").Value = "'154.63"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +0.46%  "

$ws.Range("E27").Value = "  +1.20%  "

$ws.Range("D28").Value = "'6.266"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +8.36%  "

$ws.Range("D29").Value = "'2.103"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -1.47%  "

$ws.Range("D30").Value = "'119.67"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +1.81%  "

$ws.Range("D31").Value = "'1.031"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -2.98%  "

$ws.Range("D32").Value = "'0.09564"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +0.22%  "

$ws.Range("D33").Value = "'5.548"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +2.92%  "

$ws.Range("E34").Value = "  -1.99%  "

$ws.Range("D35").Value = "'3.550"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -0.60%  "

$ws.Range("D36").Value = "'0.02267"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +0.41%  "

$ws.Range("D37").Value = "'0.06110"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +0.20%  "

$ws.Range("D38").Value = "'1.172"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -0.22%  "

$ws.Range("D39").Value = "'0.5952"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +1.10%  "

$ws.Range("D40").Value = "'10.70"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +5.53%  "

$ws.Range("D41").Value = "'7.920"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -5.33%  "

$ws.Range("D42").Value = "'0.1855"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +0.78%  "

$ws.Range("D43").Value = "'2.477"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +2.94%  "

$ws.Range("E44").Value = "  -1.03%  "

$ws.Range("D45").Value = "'0.07719"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -0.80%  "

$ws.Range("D46").Value = "'12.37"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +1.68%  "

$ws.Range("D47").Value = "'0.5573"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +0.40%  "

$ws.Range("E48").Value = "  +1.16%  "

$ws.Range("D49").Value = "'115.00"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +1.47%  "

$ws.Range("D50").Value = "'72.71"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +1.54%  "

$ws.Range("E51").Value = "  +1.96%  "
